$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# Before state (paragraphs, 1-based):
#  1 Indledning – hvorfor er emnet relevant
#  2 (empty)
#  3 Hvad er et skoleskema?
#  4 Hvad er kriterierne for et skoleskema fra ministeriet/kommunen?
#  5 Hvem er interessenterne?
#  6 Hvilke krav stilles der af interessenterne?
#  7 State of the art? Nuværende programmer.
#  8 Hvad er genetiske algoritmer?
#  9 Forskellige teorier indenfor genetiske algoritmer?
# 10 Hvordan lærer elever bedst mht skoleskema?
# 11 (empty)
# 12 (empty, contains the _GoBack bookmark)
#
# After state:
#  1 Indledning – hvorfor er emnet relevant
#  2 (empty)
#  3 Hvad er et skoleskema?
#  4 Hvad er kriterierne for et skoleskema fra ministeriet/kommunen?
#  5 Hvem er interessenterne?
#  6 Hvilke krav stilles der af interessenterne?
#  7 Hvordan lærer elever bedst med hensyn til [_GoBack]skoleskema?
#  8 Hvad er et godt skoleskema?
#  9 State of the art? Nuværende programmer.
# 10 Hvad er genetiske algoritmer?
# 11 Forskellige teorier indenfor genetiske algoritmer?
# 12 (empty)
# 13 (empty)
# ---------------------------------------------------------------------------

# Step 1: make room for the two new questions right after paragraph 6
# ("Hvilke krav stilles der af interessenterne?") by inserting two fresh
# empty paragraphs there.
$p6 = $d.Paragraphs.Item(6)
$p6.Range.InsertParagraphAfter()
$p6 = $d.Paragraphs.Item(6)
$p6.Range.InsertParagraphAfter()

# Step 2: fill the first new paragraph with the reworded question (the
# bookmark gets re-inserted into it afterwards).
$p7 = $d.Paragraphs.Item(7)
$p7.Range.Text = "Hvordan lærer elever bedst med hensyn til skoleskema? "

# Step 3: fill the second new paragraph with the brand-new question.
$p8 = $d.Paragraphs.Item(8)
$p8.Range.Text = "Hvad er et godt skoleskema?"

# Step 4: relocate the "_GoBack" bookmark so that it sits right after
# "med hensyn til " and before "skoleskema?" inside paragraph 7 (adding a
# bookmark with an existing name simply moves it, since names are unique).
$p7 = $d.Paragraphs.Item(7)
$marker = "med hensyn til "
$idx = $p7.Range.Text.IndexOf($marker)
$insertPos = $p7.Range.Start + $idx + $marker.Length
$bmRange = $d.Range($insertPos, $insertPos)
$d.Bookmarks.Add("_GoBack", $bmRange) | Out-Null

# Step 5: drop the old "Hvordan lærer elever bedst mht skoleskema?"
# paragraph now that its (reworded) content lives earlier in the document.
# It shifted down by the two paragraphs inserted in step 1, so it is now
# paragraph 12; removing it leaves the trailing two empty paragraphs intact.
$oldQuestion = $d.Paragraphs.Item(12)
$oldQuestion.Range.Delete() | Out-Null
